# Apply the benchmark-stats refresh to the single-column results table.
# The table keeps 46 rows throughout; only the text content of specific
# cells changes (rows 1-12 get new/placeholder values, and the three
# multi-value "raw log line" rows at the bottom collapse down into the
# simple summary values that used to live in rows 1-3).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row -> new cell text (1-indexed, matching Word's Cell()/Rows() indexing)
$updates = [ordered]@{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "1639"
    5  = "0.00002"
    6  = "0.00061"
    7  = "0.00015"
    8  = "0.00004"
    9  = "0.00026"
    10 = "0.00030"
    11 = "0.00037"
    12 = "0.28890"
    44 = "99.76"
    45 = "0.29"
    46 = "120"
}

foreach ($rowIndex in $updates.Keys) {
    $t.Cell($rowIndex, 1).Range.Text = $updates[$rowIndex]
}

Write-Output "Done. Row count: $($t.Rows.Count)"
